$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - (Intercept)
$ws.Range("B2").Value = 60368.859649
$ws.Range("D2").Value = 91.814448

# Row 3 - household_group_collapsed
$ws.Range("B3").Value = 17251.628704
$ws.Range("D3").Value = 13.118922
$ws.Range("E3").Value = 0.000003

# Row 4 - Residuals
$ws.Range("B4").Value = 218293.109073
$ws.Range("C4").Value = 332

# Row 5 - SM-Control
$ws.Range("G5").Value = -3.131749
$ws.Range("H5").Value = -12.428582
$ws.Range("I5").Value = 6.165084
$ws.Range("J5").Value = 0.707555

# Row 6 - SM + Traps-Control
$ws.Range("G6").Value = -17.217715
$ws.Range("H6").Value = -26.982542
$ws.Range("I6").Value = -7.452889
$ws.Range("J6").Value = 0.000124

# Row 7 - SM + Traps-SM
$ws.Range("G7").Value = -14.085966
$ws.Range("H7").Value = -21.428464
$ws.Range("I7").Value = -6.743469
$ws.Range("J7").Value = 0.000026
